$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the SRC/DST matrix: knot<->input line and input line<->knot entries
# were "false" but should reflect that a knot can actually connect via a
# fork/synapse (IoConnector direction/rotation bug).
$ws.Range("H5").Value = "fork"
$ws.Range("I5").Value = "synapse"
$ws.Range("D9").Value = "fork"
$ws.Range("D10").Value = "synapse"

# Highlight the corrected cells in yellow to flag the remaining smaller bug.
$ws.Range("I5").Interior.Color = 65535
$ws.Range("I5").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)

# Update the active selection to match the author's last-saved cursor position.
$ws.Range("L9").Select()
